$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.301.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.18%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.961.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.42%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "381.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.51%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.09%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.540"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.42%  "

# Row 8
$ws.Range("E8").Value = "  +0.16%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.589"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.95%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.18%  "

# Row 11
$ws.Range("E11").Value = "  -0.02%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0840"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.68%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.442.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.76%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.18%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.72%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.970.27"
$ws.Range("D16").Style = "Normal"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.987"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.64%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.266.67"
$ws.Range("D18").Style = "Normal"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.03%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.11%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.30%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0954"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.95%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.57%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.01%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.75%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.54%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.17%  "

# Row 28: 'Hedera' -> 'LEO'
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.36%  "

# Row 29: 'LEO' -> 'Hedera'
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.113"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +11.32%  "

# Row 30: 'Kaspa' -> 'Dai'
$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.08%  "

# Row 31: 'Dai' -> 'Kaspa'
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.166"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.97%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.33%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.16%  "

# Row 34: 'InjectiveProtocol' -> 'OKB'
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "50.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.46%  "

# Row 35: 'OKB' -> 'InjectiveProtocol'
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "34.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.28%  "

# Row 36
$ws.Range("E36").Value = "  -2.64%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0445"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.76%  "

# Row 38
$ws.Range("E38").Value = "  -0.05%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.30%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.98"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.44%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.87%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.115"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.52%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.52%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.24%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.49%  "

# Row 46
$ws.Range("E46").Value = "  +0.12%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.274"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.94%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.019.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.17%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.69%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0335"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.71%  "
